$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 1.609586333333333
$ws.Range("H2").Value2 = 4.828759
$ws.Range("I2").Value2 = 0.05107819292772156
$ws.Range("J2").Value2 = 0.05107819292772156
$ws.Range("M2").Value2 = 7.106976666666665
$ws.Range("N2").Value2 = 21.32093
$ws.Range("O2").Value2 = 0.1598176868560746
$ws.Range("P2").Value2 = 0.1598176868560746
$ws.Range("Q2").Value2 = 11.43929251398555
$ws.Range("R2").Value2 = 102.95363262587
$ws.Range("S2").Value2 = 0.008163198642496768
$ws.Range("T2").Value2 = 0.008163198642496766

# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1.609586333333333
$ws.Range("H3").Value2 = 4.828759
$ws.Range("I3").Value2 = 0.05107819292772156
$ws.Range("J3").Value2 = 0.05107819292772156
$ws.Range("O3").Value2 = 0.6350325402576649
$ws.Range("P3").Value2 = 0.6350325402576648
$ws.Range("Q3").Value2 = 45.45381131969889
$ws.Range("R3").Value2 = 409.08430187729
$ws.Range("S3").Value2 = 0.03243631460666212
$ws.Range("T3").Value2 = 0.03243631460666211

# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 1.609586333333333
$ws.Range("H4").Value2 = 4.828759
$ws.Range("I4").Value2 = 0.05107819292772156
$ws.Range("J4").Value2 = 0.05107819292772156
$ws.Range("O4").Value2 = 0.2051497728862606
$ws.Range("P4").Value2 = 0.2051497728862606
$ws.Range("Q4").Value2 = 14.68403345955722
$ws.Range("R4").Value2 = 132.156301136015
$ws.Range("S4").Value2 = 0.01047867967856268
$ws.Range("T4").Value2 = 0.01047867967856268

# Row 5
$ws.Range("I5").Value2 = 0.5992082897496871
$ws.Range("J5").Value2 = 0.5992082897496871
$ws.Range("M5").Value2 = 7.106976666666665
$ws.Range("N5").Value2 = 21.32093
$ws.Range("O5").Value2 = 0.1598176868560746
$ws.Range("P5").Value2 = 0.1598176868560746
$ws.Range("Q5").Value2 = 134.1965819533044
$ws.Range("R5").Value2 = 1207.76923757974
$ws.Range("S5").Value2 = 0.0957640828127795
$ws.Range("T5").Value2 = 0.09576408281277948

# Row 6
$ws.Range("I6").Value2 = 0.5992082897496871
$ws.Range("J6").Value2 = 0.5992082897496871
$ws.Range("O6").Value2 = 0.6350325402576649
$ws.Range("P6").Value2 = 0.6350325402576648
$ws.Range("S6").Value2 = 0.3805167623831947
$ws.Range("T6").Value2 = 0.3805167623831946

# Row 7
$ws.Range("I7").Value2 = 0.5992082897496871
$ws.Range("J7").Value2 = 0.5992082897496871
$ws.Range("O7").Value2 = 0.2051497728862606
$ws.Range("P7").Value2 = 0.2051497728862606
$ws.Range("S7").Value2 = 0.122927444553713
$ws.Range("T7").Value2 = 0.1229274445537129

# Row 8
$ws.Range("H8").Value2 = 33.060729
$ws.Range("I8").Value2 = 0.3497135173225914
$ws.Range("J8").Value2 = 0.3497135173225914
$ws.Range("M8").Value2 = 7.106976666666665
$ws.Range("N8").Value2 = 21.32093
$ws.Range("O8").Value2 = 0.1598176868560746
$ws.Range("P8").Value2 = 0.1598176868560746
$ws.Range("Q8").Value2 = 78.32060986199666
$ws.Range("R8").Value2 = 704.8854887579699
$ws.Range("S8").Value2 = 0.05589040540079833
$ws.Range("T8").Value2 = 0.05589040540079832

# Row 9
$ws.Range("H9").Value2 = 33.060729
$ws.Range("I9").Value2 = 0.3497135173225914
$ws.Range("J9").Value2 = 0.3497135173225914
$ws.Range("O9").Value2 = 0.6350325402576649
$ws.Range("P9").Value2 = 0.6350325402576648
$ws.Range("Q9").Value2 = 311.2054542497767
$ws.Range("S9").Value2 = 0.2220794632678081
$ws.Range("T9").Value2 = 0.2220794632678081

# Row 10
$ws.Range("H10").Value2 = 33.060729
$ws.Range("I10").Value2 = 0.3497135173225914
$ws.Range("J10").Value2 = 0.3497135173225914
$ws.Range("O10").Value2 = 0.2051497728862606
$ws.Range("P10").Value2 = 0.2051497728862606
$ws.Range("R10").Value2 = 904.8253717984651
$ws.Range("S10").Value2 = 0.071743648653985
$ws.Range("T10").Value2 = 0.071743648653985
